$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 156; this shifts rows 156:199 down to 157:200
# and extends the sheet dimension to A1:R200 automatically.
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new record's data.
$ws.Cells.Item(156,1).Value = 4
$ws.Cells.Item(156,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(156,3).Value = "Los Lagos"
$ws.Cells.Item(156,4).Value = 44551
$ws.Cells.Item(156,5).Value = 10
$ws.Cells.Item(156,6).Value = 100112043
$ws.Cells.Item(156,7).Value = "Pepino ensalada"
$ws.Cells.Item(156,8).Value = "Sin especificar"
$ws.Cells.Item(156,9).Value = "Primera"
$ws.Cells.Item(156,10).Value = 400
$ws.Cells.Item(156,11).Value = 12000
$ws.Cells.Item(156,12).Value = 12000
$ws.Cells.Item(156,13).Value = 12000
$ws.Cells.Item(156,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(156,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(156,16).Value = 200
$ws.Cells.Item(156,17).Value = 60
$ws.Cells.Item(156,18).Value = "Hortaliza"

# Ensure the date cell keeps the same number format style used by the rest
# of column D (style index 2 in the original sheet, i.e. same as D155/D157).
$ws.Cells.Item(156,4).NumberFormat = $ws.Cells.Item(157,4).NumberFormat
